# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.856.34"
$ws.Range("E2").Value = "  -0.10%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.904.92"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4
$ws.Range("E4").Value = "  +0.16%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.47"
$ws.Range("E5").Value = "  -3.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.50"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.901.87"
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("E9").Value = "  -0.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").Value = "  -1.34%  "

# Row 11
$ws.Range("E11").Value = "  -2.06%  "

# Row 12
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("E13").Value = "  -1.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.33"
$ws.Range("E14").Value = "  +0.29%  "

# Row 15
$ws.Range("E15").Value = "  +0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.387.51"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.810.82"
$ws.Range("E17").Value = "  -0.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.919.78"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -1.60%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.00"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("E21").Value = "  -2.67%  "

# Row 22
$ws.Range("E22").Value = "  -1.28%  "

# Row 23
$ws.Range("E23").Value = "  -0.63%  "

# Row 24
$ws.Range("E24").Value = "  -1.43%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("E25").Value = "  +1.05%  "

# Row 26
$ws.Range("E26").Value = "  -7.09%  "

# Row 27
$ws.Range("E27").Value = "  -0.02%  "

# Row 29
$ws.Range("E29").Value = "  +11.35%  "

# Row 30
$ws.Range("E30").Value = "  -3.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("E31").Value = "  -2.39%  "

# Row 32
$ws.Range("E32").Value = "  -4.18%  "

# Row 33
$ws.Range("E33").Value = "  +0.13%  "

# Row 34
$ws.Range("E34").Value = "  -2.30%  "

# Row 35
$ws.Range("E35").Value = "  -1.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("E36").Value = "  -2.23%  "

# Row 37
$ws.Range("E37").Value = "  -2.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.81"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.89"
$ws.Range("E39").Value = "  -6.72%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  -4.42%  "

# Row 41
$ws.Range("E41").Value = "  -1.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.05"
$ws.Range("E42").Value = "  +5.06%  "

# Row 43
$ws.Range("E43").Value = "  -2.59%  "

# Row 44
$ws.Range("E44").Value = "  -2.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.701.19"
$ws.Range("E45").Value = "  +0.46%  "

# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.68"
$ws.Range("E46").Value = "  -1.98%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0336"
$ws.Range("E47").Value = "  -0.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "347.46"
$ws.Range("E48").Value = "  -0.85%  "

# Row 49
$ws.Range("E49").Value = "  +0.01%  "

# Row 50
$ws.Range("E50").Value = "  -0.94%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.57"
$ws.Range("E51").Value = "  -3.70%  "
